# Apply cell-value updates to the Hades_Profits tables across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
# Generated from the scheduled-runner market-data refresh diff: per-row price/profit recalculations
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N), plus two cells that the refresh
# leaves blank this cycle (ARM!N102, CUL!M106) cleared via ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 3477.25
$ws.Range("J52").Value = 3900
$ws.Range("L52").Value = 11700
$ws.Range("N52").Value = -12020
$ws.Range("H62").Value = 2452.2593
$ws.Range("I62").Value = 2358.4119
$ws.Range("K62").Value = 2358.4119
$ws.Range("M62").Value = -1734.4119
$ws.Range("H64").Value = 4205.2905
$ws.Range("I64").Value = 3900
$ws.Range("J64").Value = 4330.1816
$ws.Range("K64").Value = 3900
$ws.Range("L64").Value = 4330.1816
$ws.Range("M64").Value = -3652
$ws.Range("N64").Value = -4826.1816
$ws.Range("H65").Value = 2452.2593
$ws.Range("I65").Value = 2358.4119
$ws.Range("K65").Value = 11792.0595
$ws.Range("M65").Value = -8672.059499999999
$ws.Range("H67").Value = 4205.2905
$ws.Range("I67").Value = 3900
$ws.Range("J67").Value = 4330.1816
$ws.Range("K67").Value = 3900
$ws.Range("L67").Value = 4330.1816
$ws.Range("M67").Value = -3042
$ws.Range("N67").Value = -6046.1816
$ws.Range("H74").Value = 3952.2104
$ws.Range("I74").Value = 3917.5454
$ws.Range("J74").Value = 3999.875
$ws.Range("K74").Value = 3917.5454
$ws.Range("L74").Value = 3999.875
$ws.Range("M74").Value = -2981.5454
$ws.Range("N74").Value = -5871.875
$ws.Range("H76").Value = 3524.2856
$ws.Range("I76").Value = 3372.7273
$ws.Range("J76").Value = 3691
$ws.Range("K76").Value = 3372.7273
$ws.Range("L76").Value = 3691
$ws.Range("M76").Value = -3057.7273
$ws.Range("N76").Value = -4321
$ws.Range("H77").Value = 3952.2104
$ws.Range("I77").Value = 3917.5454
$ws.Range("J77").Value = 3999.875
$ws.Range("K77").Value = 19587.727
$ws.Range("L77").Value = 19999.375
$ws.Range("M77").Value = -14907.727
$ws.Range("N77").Value = -29359.375
$ws.Range("H79").Value = 3524.2856
$ws.Range("I79").Value = 3372.7273
$ws.Range("J79").Value = 3691
$ws.Range("K79").Value = 3372.7273
$ws.Range("L79").Value = 3691
$ws.Range("M79").Value = -2280.7273
$ws.Range("N79").Value = -5875
$ws.Range("H98").Value = 1010.4
$ws.Range("I98").Value = 708.48
$ws.Range("J98").Value = 2520
$ws.Range("K98").Value = 708.48
$ws.Range("L98").Value = 2520
$ws.Range("M98").Value = 789.52
$ws.Range("N98").Value = -5516
$ws.Range("H122").Value = 1010.4
$ws.Range("I122").Value = 708.48
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 2125.44
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = 324.5599999999999
$ws.Range("N122").Value = -12460
$ws.Range("H135").Value = 47252.863
$ws.Range("I135").Value = 67931.07000000001
$ws.Range("J135").Value = 2942.4285
$ws.Range("K135").Value = 611379.6300000001
$ws.Range("L135").Value = 26481.8565
$ws.Range("M135").Value = -608844.6300000001
$ws.Range("N135").Value = -31551.8565
$ws.Range("H137").Value = 2779336.2
$ws.Range("I137").Value = 6251072.5
$ws.Range("J137").Value = 1947.1
$ws.Range("K137").Value = 18753217.5
$ws.Range("L137").Value = 5841.299999999999
$ws.Range("M137").Value = -18750667.5
$ws.Range("N137").Value = -10941.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15140.375
$ws.Range("I32").Value = 15410.467
$ws.Range("K32").Value = 15410.467
$ws.Range("M32").Value = -15123.467
$ws.Range("H63").Value = 3568.6667
$ws.Range("I63").Value = 2800
$ws.Range("J63").Value = 4337.3335
$ws.Range("K63").Value = 2800
$ws.Range("L63").Value = 4337.3335
$ws.Range("M63").Value = -2114
$ws.Range("N63").Value = -5709.3335
$ws.Range("H66").Value = 3568.6667
$ws.Range("I66").Value = 2800
$ws.Range("J66").Value = 4337.3335
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 21686.6675
$ws.Range("M66").Value = -10568
$ws.Range("N66").Value = -28550.6675
$ws.Range("H88").Value = 5431.6113
$ws.Range("I88").Value = 2646.9
$ws.Range("J88").Value = 8912.5
$ws.Range("K88").Value = 2646.9
$ws.Range("L88").Value = 8912.5
$ws.Range("M88").Value = -2240.9
$ws.Range("N88").Value = -9724.5
$ws.Range("H91").Value = 5431.6113
$ws.Range("I91").Value = 2646.9
$ws.Range("J91").Value = 8912.5
$ws.Range("K91").Value = 2646.9
$ws.Range("L91").Value = 8912.5
$ws.Range("M91").Value = -1242.9
$ws.Range("N91").Value = -11720.5
$ws.Range("H102").Value = 47620144
$ws.Range("I102").Value = 47620144
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 47620144
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -47618522
$ws.Range("N102").ClearContents()
$ws.Range("H138").Value = 40328.824
$ws.Range("J138").Value = 40328.824
$ws.Range("L138").Value = 40328.824
$ws.Range("N138").Value = -50608.824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1269.6666
$ws.Range("I16").Value = 900.2
$ws.Range("J16").Value = 1533.5714
$ws.Range("K16").Value = 900.2
$ws.Range("L16").Value = 1533.5714
$ws.Range("M16").Value = -613.2
$ws.Range("N16").Value = -2107.5714
$ws.Range("H20").Value = 50749.75
$ws.Range("J20").Value = 50749.75
$ws.Range("L20").Value = 50749.75
$ws.Range("N20").Value = -51221.75
$ws.Range("H30").Value = 50749.75
$ws.Range("J30").Value = 50749.75
$ws.Range("L30").Value = 50749.75
$ws.Range("N30").Value = -50931.75
$ws.Range("H31").Value = 3089.7273
$ws.Range("I31").Value = 1152.3462
$ws.Range("J31").Value = 10285.714
$ws.Range("K31").Value = 1152.3462
$ws.Range("L31").Value = 10285.714
$ws.Range("M31").Value = -857.3462
$ws.Range("N31").Value = -10875.714
$ws.Range("H34").Value = 3089.7273
$ws.Range("I34").Value = 1152.3462
$ws.Range("J34").Value = 10285.714
$ws.Range("K34").Value = 1152.3462
$ws.Range("L34").Value = 10285.714
$ws.Range("M34").Value = -950.3462
$ws.Range("N34").Value = -10689.714
$ws.Range("H62").Value = 3099.2
$ws.Range("I62").Value = 2997.4
$ws.Range("J62").Value = 3201
$ws.Range("K62").Value = 2997.4
$ws.Range("L62").Value = 3201
$ws.Range("M62").Value = -2373.4
$ws.Range("N62").Value = -4449
$ws.Range("H65").Value = 3099.2
$ws.Range("I65").Value = 2997.4
$ws.Range("J65").Value = 3201
$ws.Range("K65").Value = 14987
$ws.Range("L65").Value = 16005
$ws.Range("M65").Value = -11867
$ws.Range("N65").Value = -22245
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H113").Value = 1269.6666
$ws.Range("I113").Value = 900.2
$ws.Range("J113").Value = 1533.5714
$ws.Range("K113").Value = 900.2
$ws.Range("L113").Value = 1533.5714
$ws.Range("M113").Value = 1269.8
$ws.Range("N113").Value = -5873.5714
$ws.Range("H122").Value = 1403.8148
$ws.Range("I122").Value = 776.3333
$ws.Range("K122").Value = 2328.9999
$ws.Range("M122").Value = 121.0001000000002
$ws.Range("H128").Value = 50749.75
$ws.Range("J128").Value = 50749.75
$ws.Range("L128").Value = 50749.75
$ws.Range("N128").Value = -60709.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1130
$ws.Range("I5").Value = 1186.6666
$ws.Range("K5").Value = 3559.9998
$ws.Range("M5").Value = -3447.9998
$ws.Range("H12").Value = 516.48
$ws.Range("I12").Value = 1155.3334
$ws.Range("J12").Value = 157.125
$ws.Range("K12").Value = 3466.0002
$ws.Range("L12").Value = 471.375
$ws.Range("M12").Value = -3293.0002
$ws.Range("N12").Value = -817.375
$ws.Range("H68").Value = 1117.6923
$ws.Range("I68").Value = 734.13794
$ws.Range("J68").Value = 2230
$ws.Range("K68").Value = 2202.41382
$ws.Range("L68").Value = 6690
$ws.Range("M68").Value = -1391.41382
$ws.Range("N68").Value = -8312
$ws.Range("H71").Value = 1117.6923
$ws.Range("I71").Value = 734.13794
$ws.Range("J71").Value = 2230
$ws.Range("K71").Value = 6607.241459999999
$ws.Range("L71").Value = 20070
$ws.Range("M71").Value = -2551.241459999999
$ws.Range("N71").Value = -28182
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 10500
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -12392
$ws.Range("H132").Value = 1572.6316
$ws.Range("J132").Value = 1884.6154
$ws.Range("L132").Value = 16961.5386
$ws.Range("N132").Value = -22021.5386
$ws.Range("H135").Value = 1130
$ws.Range("I135").Value = 1186.6666
$ws.Range("K135").Value = 10679.9994
$ws.Range("M135").Value = -8144.999400000001
$ws.Range("H140").Value = 2232.95
$ws.Range("I140").Value = 2076.5908
$ws.Range("K140").Value = 6229.7724
$ws.Range("M140").Value = -1049.7724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 82076.16
$ws.Range("I70").Value = 129548.75
$ws.Range("J70").Value = 6120
$ws.Range("K70").Value = 129548.75
$ws.Range("L70").Value = 6120
$ws.Range("M70").Value = -129278.75
$ws.Range("N70").Value = -6660
$ws.Range("H73").Value = 82076.16
$ws.Range("I73").Value = 129548.75
$ws.Range("J73").Value = 6120
$ws.Range("K73").Value = 129548.75
$ws.Range("L73").Value = 6120
$ws.Range("M73").Value = -128612.75
$ws.Range("N73").Value = -7992
$ws.Range("H80").Value = 3676.9092
$ws.Range("I80").Value = 3032.2222
$ws.Range("J80").Value = 4123.231
$ws.Range("K80").Value = 3032.2222
$ws.Range("L80").Value = 4123.231
$ws.Range("M80").Value = -2034.2222
$ws.Range("N80").Value = -6119.231
$ws.Range("H83").Value = 3676.9092
$ws.Range("I83").Value = 3032.2222
$ws.Range("J83").Value = 4123.231
$ws.Range("K83").Value = 15161.111
$ws.Range("L83").Value = 20616.155
$ws.Range("M83").Value = -10169.111
$ws.Range("N83").Value = -30600.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 130.4
$ws.Range("I55").Value = 62.5
$ws.Range("J55").Value = 402
$ws.Range("K55").Value = 62.5
$ws.Range("L55").Value = 402
$ws.Range("M55").Value = 110.5
$ws.Range("N55").Value = -748
$ws.Range("H132").Value = 115408.664
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 171663
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 514989
$ws.Range("M132").Value = -6170
$ws.Range("N132").Value = -520049

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 390.75
$ws.Range("I107").Value = 615
$ws.Range("J107").Value = 316
$ws.Range("K107").Value = 1845
$ws.Range("L107").Value = 948
$ws.Range("M107").Value = 75
$ws.Range("N107").Value = -4788
$ws.Range("H113").Value = 796.4194
$ws.Range("I113").Value = 617.1053000000001
$ws.Range("J113").Value = 1080.3334
$ws.Range("K113").Value = 1851.3159
$ws.Range("L113").Value = 3241.0002
$ws.Range("M113").Value = 318.6840999999999
$ws.Range("N113").Value = -7581.0002

